$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Skew" header in F1 and the SKEW.P formula in F2
$ws.Range("F1").Value = "Skew"
$ws.Range("F2").Formula = "=SKEW.P(B2:B11)"

# Update the active selection to reflect where the new formula was entered
$ws.Range("F2").Select()
